$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "180-56316632"
$ws.Range("B2").Value = "DJAUSA4235705"
$ws.Range("C2").Value = "24N0016724"

# Row 3
$ws.Range("A3").Value = "180-57374682"
$ws.Range("B3").Value = "DJDENA4235532"
$ws.Range("C3").Value = 2160192291

# Row 4
$ws.Range("A4").Value = "180-46207486"
$ws.Range("B4").Value = 2711778683
$ws.Range("C4").Value = 2711778683

# Row 5
$ws.Range("A5").Value = "180-57261886"
$ws.Range("B5").Value = "DJDFWA4235268"
$ws.Range("C5").Value = 2233472917

# Row 6
$ws.Range("A6").Value = "180-57261923"
$ws.Range("B6").Value = "T030504439"
$ws.Range("C6").Value = "T030504439"

# Row 7
$ws.Range("A7").Value = "180-59125463"
$ws.Range("B7").Value = "DJJFKA4235203"
$ws.Range("C7").Value = 2061828984

# Row 8
$ws.Range("A8").Value = "180-61148452"
$ws.Range("B8").Value = "T070216369"
$ws.Range("C8").Value = "T070216369"

# Row 9
$ws.Range("A9").Value = "180-57261875"
$ws.Range("B9").Value = "DJJFKA4234415"
$ws.Range("C9").Value = 2061842612

# Row 10
$ws.Range("A10").Value = "180-56316444"
$ws.Range("B10").Value = "DJAUSA4234355"
$ws.Range("C10").Value = "24N0016667"

# Row 11
$ws.Range("A11").Value = "180-56316470"
$ws.Range("B11").Value = "DJJFKA4233952"
$ws.Range("C11").Value = 2061844205
